$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:O25 (row-major), matching the updated case data ("case with 380 kV done").
$data = New-Object 'object[,]' 24,14
$data[0,0] = 0.2882132122826988
$data[0,1] = 0.04006109076328812
$data[0,2] = 0.2543585930202283
$data[0,3] = 0
$data[0,4] = 1.556594914100124
$data[0,5] = 0.00247158008162638
$data[0,6] = 0
$data[0,7] = 0.768751709634742
$data[0,8] = 0.3073440924463711
$data[0,9] = 0.3037661563881215
$data[0,10] = 0
$data[0,11] = 0.2730817823166234
$data[0,12] = 0
$data[0,13] = 3.62665395167744
$data[1,0] = 0.2563337833487935
$data[1,1] = 0.03504299120776011
$data[1,2] = 0.2490823709131149
$data[1,3] = 0
$data[1,4] = 1.562082103744899
$data[1,5] = 0.002473959002247813
$data[1,6] = 0
$data[1,7] = 0.7759094548576009
$data[1,8] = 0.3058834325174828
$data[1,9] = 0.2689076935155867
$data[1,10] = 0
$data[1,11] = 0.2595079266723488
$data[1,12] = 0
$data[1,13] = 3.650691019819192
$data[2,0] = 0.2367600123551483
$data[2,1] = 0.031953788277022
$data[2,2] = 0.2459465689774305
$data[2,3] = 0
$data[2,4] = 1.566232761400791
$data[2,5] = 0.002475498054288273
$data[2,6] = 0
$data[2,7] = 0.7806806984437777
$data[2,8] = 0.3051430898805805
$data[2,9] = 0.2474913289713214
$data[2,10] = 0
$data[2,11] = 0.2512761700026616
$data[2,12] = 0
$data[2,13] = 3.66729933361367
$data[3,0] = 0.2287841212755097
$data[3,1] = 0.0306929443917312
$data[3,2] = 0.2446949278645434
$data[3,3] = 0
$data[3,4] = 1.568120812025718
$data[3,5] = 0.002476144998871488
$data[3,6] = 0
$data[3,7] = 0.7827196588047549
$data[3,8] = 0.3048807777764395
$data[3,9] = 0.2387612012149418
$data[3,10] = 0
$data[3,11] = 0.2479476766143662
$data[3,12] = 0
$data[3,13] = 3.674532481536929
$data[4,0] = 0.227459780311932
$data[4,1] = 0.0304834651157222
$data[4,2] = 0.2444886811089333
$data[4,3] = 0
$data[4,4] = 1.568446200211483
$data[4,5] = 0.002476253619135447
$data[4,6] = 0
$data[4,7] = 0.7830639426766304
$data[4,8] = 0.3048396005534784
$data[4,9] = 0.2373114171440704
$data[4,10] = 0
$data[4,11] = 0.2473965600137191
$data[4,12] = 0
$data[4,13] = 3.675761634569881
$data[5,0] = 0.2366524436915824
$data[5,1] = 0.03193679198660959
$data[5,2] = 0.245929582573325
$data[5,3] = 0
$data[5,4] = 1.566257428045411
$data[5,5] = 0.002475506698962904
$data[5,6] = 0
$data[5,7] = 0.7807078133820582
$data[5,8] = 0.3051393927477974
$data[5,9] = 0.2473736020081958
$data[5,10] = 0
$data[5,11] = 0.2512311751758958
$data[5,12] = 0
$data[5,13] = 3.667394999020601
$data[6,0] = 0.2772214052493496
$data[6,1] = 0.03833255960796578
$data[6,2] = 0.2525178720542129
$data[6,3] = 0
$data[6,4] = 1.558324759679593
$data[6,5] = 0.0024723840993569
$data[6,6] = 0
$data[6,7] = 0.7711415878168957
$data[6,8] = 0.3068079892723574
$data[6,9] = 0.291749981168806
$data[6,10] = 0
$data[6,11] = 0.2683803299744767
$data[6,12] = 0
$data[6,13] = 3.634558114205291
$data[7,0] = 0.3567616941825804
$data[7,1] = 0.05080855616512281
$data[7,2] = 0.2662569771014205
$data[7,3] = 0
$data[7,4] = 1.548966081653916
$data[7,5] = 0.002466879999272103
$data[7,6] = 0
$data[7,7] = 0.7553684716751761
$data[7,8] = 0.3113212405861745
$data[7,9] = 0.378649376119256
$data[7,10] = 0
$data[7,11] = 0.3028171862646403
$data[7,12] = 0
$data[7,13] = 3.584838411757715
$data[8,0] = 0.4151723077485485
$data[8,1] = 0.05993250993424226
$data[8,2] = 0.2768462338728455
$data[8,3] = 0
$data[8,4] = 1.545865212220988
$data[8,5] = 0.002463209977269934
$data[8,6] = 0
$data[8,7] = 0.7456004849721936
$data[8,8] = 0.3153936555668935
$data[8,9] = 0.4424012219393774
$data[8,10] = 0
$data[8,11] = 0.3286035171193191
$data[8,12] = 0
$data[8,13] = 3.557255179384185
$data[9,0] = 0.4417352967986687
$data[9,1] = 0.06407373759810753
$data[9,2] = 0.281770211904572
$data[9,3] = 0
$data[9,4] = 1.545273719586589
$data[9,5] = 0.002461620777164136
$data[9,6] = 0
$data[9,7] = 0.7415521114058627
$data[9,8] = 0.3174106068478153
$data[9,9] = 0.4713798554942343
$data[9,10] = 0
$data[9,11] = 0.3404385862284371
$data[9,12] = 0
$data[9,13] = 3.546649682411868
$data[10,0] = 0.4517923935321733
$data[10,1] = 0.06564052609540738
$data[10,2] = 0.2836500615634492
$data[10,3] = 0
$data[10,4] = 1.545167450950686
$data[10,5] = 0.002461030477548332
$data[10,6] = 0
$data[10,7] = 0.7400759250156597
$data[10,8] = 0.318197998320727
$data[10,9] = 0.482349653598078
$data[10,10] = 0
$data[10,11] = 0.3449351131577885
$data[10,12] = 0
$data[10,13] = 3.542912955648603
$data[11,0] = 0.4496265033344287
$data[11,1] = 0.06530315368399897
$data[11,2] = 0.2832445261759915
$data[11,3] = 0
$data[11,4] = 1.545185103438172
$data[11,5] = 0.002461157098766393
$data[11,6] = 0
$data[11,7] = 0.7403913195769931
$data[11,8] = 0.3180273695315634
$data[11,9] = 0.4799872883876901
$data[11,10] = 0
$data[11,11] = 0.3439660489505272
$data[11,12] = 0
$data[11,13] = 3.543705301877566
$data[12,0] = 0.4425627379633283
$data[12,1] = 0.06420266675965536
$data[12,2] = 0.281924563412403
$data[12,3] = 0
$data[12,4] = 1.545262617985387
$data[12,5] = 0.002461571982741403
$data[12,6] = 0
$data[12,7] = 0.7414295252638681
$data[12,8] = 0.3174749128564258
$data[12,9] = 0.4722824263764949
$data[12,10] = 0
$data[12,11] = 0.3408082219758768
$data[12,12] = 0
$data[12,13] = 3.546336660877699
$data[13,0] = 0.4382357363712117
$data[13,1] = 0.06352840185950015
$data[13,2] = 0.28111802996915
$data[13,3] = 0
$data[13,4] = 1.545325425996396
$data[13,5] = 0.002461827606824005
$data[13,6] = 0
$data[13,7] = 0.7420728598465942
$data[13,8] = 0.3171395919724915
$data[13,9] = 0.4675624668588512
$data[13,10] = 0
$data[13,11] = 0.3388758883140497
$data[13,12] = 0
$data[13,13] = 3.547984824680441
$data[14,0] = 0.4134361410297913
$data[14,1] = 0.0596616780389212
$data[14,2] = 0.2765265804369932
$data[14,3] = 0
$data[14,4] = 1.545920344240599
$data[14,5] = 0.002463315446900011
$data[14,6] = 0
$data[14,7] = 0.7458730089407446
$data[14,8] = 0.3152651490131433
$data[14,9] = 0.440506901387522
$data[14,10] = 0
$data[14,11] = 0.3278321564707625
$data[14,12] = 0
$data[14,13] = 3.557987361045718
$data[15,0] = 0.3982198851832379
$data[15,1] = 0.05728713764061411
$data[15,2] = 0.2737371605431065
$data[15,3] = 0
$data[15,4] = 1.546495051011021
$data[15,5] = 0.002464248720806292
$data[15,6] = 0
$data[15,7] = 0.7483054982158563
$data[15,8] = 0.3141573265871642
$data[15,9] = 0.4239030516888818
$data[15,10] = 0
$data[15,11] = 0.3210838553914428
$data[15,12] = 0
$data[15,13] = 3.564621087477605
$data[16,0] = 0.3894671488631047
$data[16,1] = 0.05592049297123936
$data[16,2] = 0.2721428265962231
$data[16,3] = 0
$data[16,4] = 1.546902709724563
$data[16,5] = 0.002464793077488568
$data[16,6] = 0
$data[16,7] = 0.749741791947379
$data[16,8] = 0.3135356109272038
$data[16,9] = 0.4143508814021573
$data[16,10] = 0
$data[16,11] = 0.3172122825524468
$data[16,12] = 0
$data[16,13] = 3.568619442565193
$data[17,0] = 0.3865035106292396
$data[17,1] = 0.05545762286237732
$data[17,2] = 0.2716047451809089
$data[17,3] = 0
$data[17,4] = 1.54705398052883
$data[17,5] = 0.002464978687771968
$data[17,6] = 0
$data[17,7] = 0.7502344830176035
$data[17,8] = 0.3133277671498718
$data[17,9] = 0.4111163413730026
$data[17,10] = 0
$data[17,11] = 0.3159031365894407
$data[17,12] = 0
$data[17,13] = 3.570004613319355
$data[18,0] = 0.3998397629573844
$data[18,1] = 0.05754000223569733
$data[18,2] = 0.2740330580819403
$data[18,3] = 0
$data[18,4] = 1.546425893170436
$data[18,5] = 0.002464148589785101
$data[18,6] = 0
$data[18,7] = 0.7480427063460802
$data[18,8] = 0.3142736547903553
$data[18,9] = 0.4256707786636014
$data[18,10] = 0
$data[18,11] = 0.3218012034833038
$data[18,12] = 0
$data[18,13] = 3.563895995334605
$data[19,0] = 0.4446375883350981
$data[19,1] = 0.06452594509352139
$data[19,2] = 0.2823118556716793
$data[19,3] = 0
$data[19,4] = 1.545236655845017
$data[19,5] = 0.002461449809611994
$data[19,6] = 0
$data[19,7] = 0.7411230361340593
$data[19,8] = 0.3176365420337532
$data[19,9] = 0.4745456382036082
$data[19,10] = 0
$data[19,11] = 0.3417353512050099
$data[19,12] = 0
$data[19,13] = 3.545556185214423
$data[20,0] = 0.4739052202308471
$data[20,1] = 0.06908344514049247
$data[20,2] = 0.2878113374308242
$data[20,3] = 0
$data[20,4] = 1.545145531122557
$data[20,5] = 0.00245975298458935
$data[20,6] = 0
$data[20,7] = 0.7369319784417598
$data[20,8] = 0.3199720232606182
$data[20,9] = 0.5064658839794731
$data[20,10] = 0
$data[20,11] = 0.3548498849590729
$data[20,12] = 0
$data[20,13] = 3.535198317362074
$data[21,0] = 0.4582856575215999
$data[21,1] = 0.06665179541168698
$data[21,2] = 0.2848680748184051
$data[21,3] = 0
$data[21,4] = 1.545131410221288
$data[21,5] = 0.002460652500228509
$data[21,6] = 0
$data[21,7] = 0.7391384982203633
$data[21,8] = 0.3187129466236058
$data[21,9] = 0.4894316732242885
$data[21,10] = 0
$data[21,11] = 0.3478425764556263
$data[21,12] = 0
$data[21,13] = 3.540577503971633
$data[22,0] = 0.3991074308410418
$data[22,1] = 0.05742568678770965
$data[22,2] = 0.2738992537012024
$data[22,3] = 0
$data[22,4] = 1.546456918788678
$data[22,5] = 0.002464193834819332
$data[22,6] = 0
$data[22,7] = 0.7481613967343392
$data[22,8] = 0.3142210155055096
$data[22,9] = 0.4248716091793483
$data[22,10] = 0
$data[22,11] = 0.3214768651170843
$data[22,12] = 0
$data[22,13] = 3.564223234668987
$data[23,0] = 0.3352474772067353
$data[23,1] = 0.0474407380979045
$data[23,2] = 0.2624529128588193
$data[23,3] = 0
$data[23,4] = 1.550834708140798
$data[23,5] = 0.002468303087681534
$data[23,6] = 0
$data[23,7] = 0.7593158179970665
$data[23,8] = 0.3099674240017691
$data[23,9] = 0.3551558196492977
$data[23,10] = 0
$data[23,11] = 0.2934153528473331
$data[23,12] = 0
$data[23,13] = 3.596717967307058

$ws.Range("B2:O25").Value = $data
